$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so values
# like "580.09" are not auto-converted to numbers by Excel,
# matching the original inline-string storage. Style is reset
# back to Normal afterwards so no visible formatting changes.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.382.21"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "3.624.30"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "580.09"
$ws.Range("E5").Value = "  -3.27%  "
$ws.Range("D6").Value = "175.18"
$ws.Range("E6").Value = "  -3.98%  "
$ws.Range("D7").Value = "3.614.38"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -5.39%  "
$ws.Range("D11").Value = "6.84"
$ws.Range("E11").Value = "  +18.23%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "48.16"
$ws.Range("E13").Value = "  -4.25%  "
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "4.208.84"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "675.35"
$ws.Range("E16").Value = "  -5.00%  "
$ws.Range("D17").Value = "8.88"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "3.629.88"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "70.405.94"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("E21").Value = "  -4.44%  "
$ws.Range("D22").Value = "11.44"
$ws.Range("E22").Value = "  -2.73%  "
$ws.Range("D23").Value = "0.937"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "17.05"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").Value = "99.56"
$ws.Range("E25").Value = "  -5.74%  "
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "9.84"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("D30").Value = "34.54"
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").Value = "3.28"
$ws.Range("E32").Value = "  -4.72%  "
$ws.Range("D33").Value = "7.53"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  -7.08%  "
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("D36").Value = "577.07"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").Value = "11.04"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "58.10"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "3.556.15"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").Value = "0.0450"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").Value = "34.35"
$ws.Range("E45").Value = "  -4.78%  "
$ws.Range("D46").Value = "0.0₃0727"
$ws.Range("E46").Value = "  -7.32%  "
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("D48").Value = "2.85"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("E51").Value = "  -2.09%  "

$ws.Range("D2:D51").Style = "Normal"
